$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "ZGN / Zhongshan, China" row (row 242) entirely, shifting all
# subsequent rows up by one.
$ws.Rows(242).Delete()
